# Weekly update: insert 3 new price records (newest week, 2023-03-28 = serial 45013)
# at the top of the data block (row 1039), pushing the existing rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 1039, shifting all the
# existing data (old rows 1039:1121) down to 1042:1124.
$ws.Rows("1039:1041").Insert()

# --- New row 1039 ---
$ws.Cells.Item(1039, 1).Value = 6
$ws.Cells.Item(1039, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1039, 3).Value = "Metropolitana"
$ws.Cells.Item(1039, 4).Value = 45013
$ws.Cells.Item(1039, 5).Value = 13
$ws.Cells.Item(1039, 6).Value = 100112008
$ws.Cells.Item(1039, 7).Value = "Coliflor"
$ws.Cells.Item(1039, 8).Value = "Sin especificar"
$ws.Cells.Item(1039, 9).Value = "Primera"
$ws.Cells.Item(1039, 10).Value = 3700
$ws.Cells.Item(1039, 11).Value = 900
$ws.Cells.Item(1039, 12).Value = 900
$ws.Cells.Item(1039, 13).Value = 900
$ws.Cells.Item(1039, 14).Value = "`$/unidad"
$ws.Cells.Item(1039, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1039, 16).Value = 900
$ws.Cells.Item(1039, 17).Value = 1
$ws.Cells.Item(1039, 18).Value = "Hortaliza"

# --- New row 1040 ---
$ws.Cells.Item(1040, 1).Value = 6
$ws.Cells.Item(1040, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1040, 3).Value = "Metropolitana"
$ws.Cells.Item(1040, 4).Value = 45013
$ws.Cells.Item(1040, 5).Value = 13
$ws.Cells.Item(1040, 6).Value = 100112008
$ws.Cells.Item(1040, 7).Value = "Coliflor"
$ws.Cells.Item(1040, 8).Value = "Sin especificar"
$ws.Cells.Item(1040, 9).Value = "Primera"
$ws.Cells.Item(1040, 10).Value = 5600
$ws.Cells.Item(1040, 11).Value = 1000
$ws.Cells.Item(1040, 12).Value = 1100
$ws.Cells.Item(1040, 13).Value = 1046
$ws.Cells.Item(1040, 14).Value = "`$/unidad"
$ws.Cells.Item(1040, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(1040, 16).Value = 1046
$ws.Cells.Item(1040, 17).Value = 1
$ws.Cells.Item(1040, 18).Value = "Hortaliza"

# --- New row 1041 ---
$ws.Cells.Item(1041, 1).Value = 6
$ws.Cells.Item(1041, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1041, 3).Value = "Metropolitana"
$ws.Cells.Item(1041, 4).Value = 45013
$ws.Cells.Item(1041, 5).Value = 13
$ws.Cells.Item(1041, 6).Value = 100112008
$ws.Cells.Item(1041, 7).Value = "Coliflor"
$ws.Cells.Item(1041, 8).Value = "Sin especificar"
$ws.Cells.Item(1041, 9).Value = "Segunda"
$ws.Cells.Item(1041, 10).Value = 6700
$ws.Cells.Item(1041, 11).Value = 700
$ws.Cells.Item(1041, 12).Value = 800
$ws.Cells.Item(1041, 13).Value = 752
$ws.Cells.Item(1041, 14).Value = "`$/unidad"
$ws.Cells.Item(1041, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1041, 16).Value = 752
$ws.Cells.Item(1041, 17).Value = 1
$ws.Cells.Item(1041, 18).Value = "Hortaliza"
